$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title text in A1, merged across A1:K1 (header band for the report)
$ws.Range("A1").Value = "EXPENSES TRACKING SYSTEM"

$titleRange = $ws.Range("A1:K1")
$titleRange.Merge()

# Big bold centered title font
$titleRange.Font.Bold = $true
$titleRange.Font.Size = 26
$titleRange.HorizontalAlignment = -4108

# Taller header row to fit the larger font
$ws.Rows.Item(1).RowHeight = 31.5

# Leave the selection where the sheet was last left (mirrors the saved view state)
$ws.Range("I12:I14").Select() | Out-Null
